$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The header label in D1 was shortened from "1-rho" to "rho".
$ws.Range("D1").Value = "rho"

# Reflect the resulting selection/active cell on the sheet (D1),
# matching the saved view state after the edit.
$null = $ws.Range("D1").Select()
